# Scheduled runner update: refresh cached market/profit figures on the
# Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 9321.637000000001
$ws.Range("I76").Value = 13509.1
$ws.Range("J76").Value = 5832.0835
$ws.Range("K76").Value = 13509.1
$ws.Range("L76").Value = 5832.0835
$ws.Range("M76").Value = -13194.1
$ws.Range("N76").Value = -6462.0835

$ws.Range("H79").Value = 9321.637000000001
$ws.Range("I79").Value = 13509.1
$ws.Range("J79").Value = 5832.0835
$ws.Range("K79").Value = 13509.1
$ws.Range("L79").Value = 5832.0835
$ws.Range("M79").Value = -12417.1
$ws.Range("N79").Value = -8016.0835

$ws.Range("H100").Value = 2625
$ws.Range("I100").Value = 1325.25
$ws.Range("K100").Value = 1325.25
$ws.Range("M100").Value = -784.25

$ws.Range("H138").Value = 25651642
$ws.Range("I138").Value = 1999
$ws.Range("J138").Value = 37051484
$ws.Range("K138").Value = 5997
$ws.Range("L138").Value = 111154452
$ws.Range("M138").Value = -857
$ws.Range("N138").Value = -111164732

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 967.625
$ws.Range("I2").Value = 963
$ws.Range("K2").Value = 963
$ws.Range("M2").Value = -850

$ws.Range("H28").Value = 20000
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 20000
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = ""
$ws.Range("M28").Value = 20000
$ws.Range("N28").Value = -20384

$ws.Range("H31").Value = 7636
$ws.Range("I31").Value = 9181.333000000001
$ws.Range("J31").Value = 3000
$ws.Range("K31").Value = 9181.333000000001
$ws.Range("L31").Value = 3000
$ws.Range("M31").Value = -8887.333000000001
$ws.Range("N31").Value = -3588

$ws.Range("H61").Value = 11368059
$ws.Range("I61").Value = 12504115
$ws.Range("K61").Value = 12504115
$ws.Range("M61").Value = -12503903

$ws.Range("H75").Value = 56333.332
$ws.Range("J75").Value = 44500
$ws.Range("L75").Value = 44500
$ws.Range("N75").Value = -46248

$ws.Range("H76").Value = 77695
$ws.Range("J76").Value = 77695
$ws.Range("L76").Value = 77695
$ws.Range("N76").Value = -78371

$ws.Range("H78").Value = 56333.332
$ws.Range("J78").Value = 44500
$ws.Range("L78").Value = 133500
$ws.Range("N78").Value = -142236

$ws.Range("H79").Value = 77695
$ws.Range("J79").Value = 77695
$ws.Range("L79").Value = 77695
$ws.Range("N79").Value = -80035

$ws.Range("H99").Value = 20000
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 20000
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = ""
$ws.Range("M99").Value = 20000
$ws.Range("N99").Value = -25990

$ws.Range("H104").Value = 20000
$ws.Range("J104").Value = 20000
$ws.Range("L104").Value = 20000
$ws.Range("N104").Value = -26988

$ws.Range("H116").Value = 967.625
$ws.Range("I116").Value = 963
$ws.Range("K116").Value = 963
$ws.Range("M116").Value = 1331

$ws.Range("H122").Value = 2830.3845
$ws.Range("I122").Value = 2106.818
$ws.Range("K122").Value = 6320.454000000001
$ws.Range("M122").Value = -3870.454000000001

$ws.Range("H132").Value = 30354318
$ws.Range("I132").Value = 10130.814
$ws.Range("K132").Value = 30392.442
$ws.Range("M132").Value = -27862.442

$ws.Range("H136").Value = 11368059
$ws.Range("I136").Value = 12504115
$ws.Range("K136").Value = 37512345
$ws.Range("M136").Value = -37509795

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 967.625
$ws.Range("I3").Value = 963
$ws.Range("K3").Value = 963
$ws.Range("M3").Value = -849

$ws.Range("H134").Value = 2585
$ws.Range("I134").Value = 2182
$ws.Range("K134").Value = 6546
$ws.Range("M134").Value = -4011

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 39593.832
$ws.Range("J68").Value = 40459
$ws.Range("L68").Value = 40459
$ws.Range("N68").Value = -41957

$ws.Range("H71").Value = 39593.832
$ws.Range("J71").Value = 40459
$ws.Range("L71").Value = 121377
$ws.Range("N71").Value = -128865

$ws.Range("H99").Value = 8976.111000000001
$ws.Range("I99").Value = 9531.200000000001
$ws.Range("K99").Value = 9531.200000000001
$ws.Range("M99").Value = -8033.200000000001

$ws.Range("H126").Value = 8976.111000000001
$ws.Range("I126").Value = 9531.200000000001
$ws.Range("K126").Value = 28593.6
$ws.Range("M126").Value = -26123.6

$ws.Range("H132").Value = 44836.168
$ws.Range("I132").Value = 52793.023
$ws.Range("J132").Value = 5051.875
$ws.Range("K132").Value = 158379.069
$ws.Range("L132").Value = 15155.625
$ws.Range("M132").Value = -155849.069
$ws.Range("N132").Value = -20215.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 593.9286
$ws.Range("I12").Value = 344
$ws.Range("J12").Value = 662.0909
$ws.Range("K12").Value = 1032
$ws.Range("L12").Value = 1986.2727
$ws.Range("M12").Value = -859
$ws.Range("N12").Value = -2332.2727

$ws.Range("H13").Value = 3336763
$ws.Range("I13").Value = 9999999
$ws.Range("J13").Value = 5145
$ws.Range("K13").Value = 29999997
$ws.Range("L13").Value = 15435
$ws.Range("M13").Value = -29999829
$ws.Range("N13").Value = -15771

$ws.Range("H112").Value = 4675
$ws.Range("I112").Value = 4650
$ws.Range("J112").Value = 4750
$ws.Range("K112").Value = 13950
$ws.Range("L112").Value = 14250
$ws.Range("M112").Value = -12842
$ws.Range("N112").Value = -16466

$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = ""
$ws.Range("M122").Value = 18000
$ws.Range("N122").Value = -22900

$ws.Range("H129").Value = 1874.5518
$ws.Range("I129").Value = 724.5
$ws.Range("J129").Value = 2312.6667
$ws.Range("K129").Value = 2173.5
$ws.Range("L129").Value = 6938.000100000001
$ws.Range("M129").Value = 2826.5
$ws.Range("N129").Value = -16938.0001

$ws.Range("H131").Value = 1630.5128
$ws.Range("I131").Value = 1118.75
$ws.Range("K131").Value = 3356.25
$ws.Range("M131").Value = 1683.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 35747.5
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 35747.5
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = ""
$ws.Range("M52").Value = 35747.5
$ws.Range("N52").Value = -36265.5

$ws.Range("H107").Value = 778.6111
$ws.Range("I107").Value = 646.8333
$ws.Range("K107").Value = 646.8333
$ws.Range("M107").Value = 1273.1667

$ws.Range("H132").Value = 3736.1875
$ws.Range("I132").Value = 3906.077
$ws.Range("K132").Value = 11718.231
$ws.Range("M132").Value = -9188.231

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3964.4614
$ws.Range("I82").Value = 2262.6667
$ws.Range("J82").Value = 4475
$ws.Range("K82").Value = 2262.6667
$ws.Range("L82").Value = 4475
$ws.Range("M82").Value = -1901.6667
$ws.Range("N82").Value = -5197

$ws.Range("H85").Value = 3964.4614
$ws.Range("I85").Value = 2262.6667
$ws.Range("J85").Value = 4475
$ws.Range("K85").Value = 2262.6667
$ws.Range("L85").Value = 4475
$ws.Range("M85").Value = -1014.6667
$ws.Range("N85").Value = -6971

$ws.Range("H100").Value = 3278.1428
$ws.Range("I100").Value = 2884.353
$ws.Range("J100").Value = 4951.75
$ws.Range("K100").Value = 2884.353
$ws.Range("L100").Value = 4951.75
$ws.Range("M100").Value = -2343.353
$ws.Range("N100").Value = -6033.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1429.8485
$ws.Range("I122").Value = 1191.6428
$ws.Range("K122").Value = 3574.9284
$ws.Range("M122").Value = -1124.9284

$ws.Range("H126").Value = 13180.667
$ws.Range("I126").Value = 16021
$ws.Range("K126").Value = 48063
$ws.Range("M126").Value = -45593
